# Update "想去人数" (F column) values for a handful of rows across sheets,
# matching the regenerated data snapshot (output generated at 7921097).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 11416
$ws1.Range("F3").Value = 1939
$ws1.Range("F5").Value = 850
$ws1.Range("F6").Value = 2426
$ws1.Range("F9").Value = 597
$ws1.Range("F11").Value = 1347
$ws1.Range("F12").Value = 678
$ws1.Range("F13").Value = 121
$ws1.Range("F14").Value = 13
$ws1.Range("F15").Value = 989
$ws1.Range("F16").Value = 538
$ws1.Range("F18").Value = 1125
$ws1.Range("F19").Value = 210
$ws1.Range("F21").Value = 8
$ws1.Range("F22").Value = 137
$ws1.Range("F23").Value = 305
$ws1.Range("F27").Value = 502
$ws1.Range("F28").Value = 677

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 37
$ws2.Range("F8").Value = 96
$ws2.Range("F10").Value = 386

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 48

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 11416
$ws4.Range("F3").Value = 1939
$ws4.Range("F6").Value = 850
$ws4.Range("F7").Value = 2426
$ws4.Range("F10").Value = 37
$ws4.Range("F11").Value = 597
$ws4.Range("F13").Value = 48
$ws4.Range("F14").Value = 1347
$ws4.Range("F16").Value = 678
$ws4.Range("F17").Value = 121
$ws4.Range("F19").Value = 13
$ws4.Range("F20").Value = 989
$ws4.Range("F21").Value = 538
$ws4.Range("F23").Value = 1125
$ws4.Range("F24").Value = 210
$ws4.Range("F26").Value = 8
$ws4.Range("F27").Value = 137
$ws4.Range("F28").Value = 305
$ws4.Range("F33").Value = 96
$ws4.Range("F34").Value = 96
$ws4.Range("F36").Value = 502
$ws4.Range("F37").Value = 677
$ws4.Range("F41").Value = 386
